$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331; existing rows 331-356 shift down to 332-357.
$ws.Rows("331:331").Insert()

# Populate the newly inserted row 331 with the new weekly price record.
$ws.Cells.Item(331, 1).Value = 10
$ws.Cells.Item(331, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(331, 3).Value = "La Araucanía"
$ws.Cells.Item(331, 4).Value = 44931
$ws.Cells.Item(331, 5).Value = 9
$ws.Cells.Item(331, 6).Value = "Fruta"
$ws.Cells.Item(331, 7).Value = 100102
$ws.Cells.Item(331, 8).Value = "Cítricos"
$ws.Cells.Item(331, 9).Value = 100102006
$ws.Cells.Item(331, 10).Value = "Pomelo"
$ws.Cells.Item(331, 11).Value = "Start Ruby"
$ws.Cells.Item(331, 12).Value = "Primera"
$ws.Cells.Item(331, 13).Value = 55
$ws.Cells.Item(331, 14).Value = 14000
$ws.Cells.Item(331, 15).Value = 14000
$ws.Cells.Item(331, 16).Value = 14000
$ws.Cells.Item(331, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(331, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(331, 19).Value = 933
$ws.Cells.Item(331, 20).Value = 15
